$wb = $excel.ActiveWorkbook

# 1. On the "PI hours" sheet, Romit Roy Choudhury's cfop list changed order
#    (RRC now listed before CHOUDHURY).
$piSheet = $wb.Worksheets.Item("PI hours")
$piSheet.Range("G3").Value = "['cfop_RRC', 'cfop_CHOUDHURY']"

# 2. Add a "users" column (E) to the "project hours" sheet, listing the
#    users who worked on each project.
$projSheet = $wb.Worksheets.Item("project hours")

# Header cell - match the bold/bordered/centered style used by the other
# header cells in row 1 (B1:D1).
$projSheet.Range("E1").Value = "users"
$projSheet.Range("B1").Copy() | Out-Null
$projSheet.Range("E1").PasteSpecial(-4122) | Out-Null

# Data rows.
$projSheet.Range("E2").Value = "['HYUNG JIN YOON', 'Mitchell Jones', 'HYUNG-JIN YOON', 'Arun Lakshmanan']"
$projSheet.Range("E3").Value = "['Ashutosh Dhekne', 'Mahanth Gowda']"
